$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update ActualRate (column E) and Result (column F) for rows 30 and 31
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "$549.37"
$ws.Range("F30").Value = "FAIL"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "$305.85"
$ws.Range("F31").Value = "FAIL"
